$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Pepino ensalada" block (row 908),
# shifting the existing rows 908-948 down to 910-950, then populate the
# two new rows with the new week's data (date 45147).
$ws.Rows("908:909").Insert()

$newRows = @(
    @{ Row = 908; D = 45147; I = "Primera"; J = 500; K = 9000;  L = 10000; M = 9500; N = "`$/caja 60 unidades"; P = 158; Q = 60 },
    @{ Row = 909; D = 45147; I = "Segunda"; J = 360; K = 6000;  L = 7000;  M = 6500; N = "`$/caja 80 unidades"; P = 81;  Q = 80 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 8
    $ws.Cells.Item($row, 2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = 100112043
    $ws.Cells.Item($row, 7).Value  = "Pepino ensalada"
    $ws.Cells.Item($row, 8).Value  = "Sin especificar"
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
